$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Update the label (A15): "h - ширина подшипника (int)" -> "B - ширина подшипника (int)"
$ws.Range("A15").Value = "B - ширина подшипника (int)"

# Update "Модель подшипника" value (B5): drop the trailing non-breaking space
$ws.Range("B5").Value = "NU 234 ECM"

# Update "Subtype" value (B4): add "roller bearings" suffix
$ws.Range("B4").Value = "Single row cilindrical roller bearings"

# Move the active selection from C12 to B5, matching the saved view state
$ws.Activate()
$ws.Range("B5").Select()
